$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algae-added")

# Copy A7:B7's formatting (date style, wrap-text bucket style) down to A8:B8
# so the new row matches the existing rows' look without inventing new styles.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row height used by the other wrapped-note rows.
$ws.Rows.Item(8).RowHeight = 32

# New feeding/count row for 2017-08-05.
$ws.Range("A8").Value = 42952
$ws.Range("B8").Value = "400 mL Ciso, 200 609, 200 Chagra"
$ws.Range("C8").Value = 195
$ws.Range("D8").Value = 215
$ws.Range("E8").Value = 290
$ws.Range("F8").Value = 216
$ws.Range("G8").Value = 238
$ws.Range("H8").Formula = "=AVERAGE(C8:G8)"
$ws.Range("I8").Formula = "=(H8*9)/0.0009"
$ws.Range("J8").Formula = "=15000*50000"
$ws.Range("K8").Formula = "=J8/I8"
$ws.Range("L8").Value = 550
$ws.Range("M8").Formula = "=L8*I8"
$ws.Range("N8").Formula = "=M8/15000"
$ws.Range("O8").Value = "250 mL Ciso, 150 609, 150 Chagra. Concentration must likely an overestimate"

# Update the saved view state (scroll/selection) to match where the editor left off.
$ws.Range("L9").Select()
